# Update "想去人数" (want-to-go count) values in the "展览" and "全部类型"
# sheets to reflect newly scraped totals.

$wb = $excel.ActiveWorkbook

# Sheet name -> { row number => new "想去人数" (column F) value }
$sheetRowMap = @{
    "展览"   = @{ 3 = 3939; 4 = 2321; 5 = 460; 8 = 191; 10 = 24; 12 = 1461; 14 = 2666 }
    "全部类型" = @{ 3 = 3939; 4 = 2321; 5 = 460; 9 = 191; 11 = 24; 15 = 1461; 17 = 2666 }
}

foreach ($sheetName in $sheetRowMap.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $sheetRowMap[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Cells.Item($row, 6).Value = $rows[$row]
    }
}
